$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 12, shifting the existing rows 12-16 down to 13-17.
$ws.Rows("12:12").Insert()

# Populate the newly inserted row 12 with the new weekly record.
$ws.Range("A12").Value = 11
$ws.Range("B12").Value = "Vega Monumental Concepción"
$ws.Range("C12").Value = "Bíobío"
$ws.Range("D12").Value = 45093
$ws.Range("E12").Value = 8
$ws.Range("F12").Value = "Fruta"
$ws.Range("G12").Value = 100107
$ws.Range("H12").Value = "Otros"
$ws.Range("I12").Value = 100107001
$ws.Range("J12").Value = "Caqui"
$ws.Range("K12").Value = "Mankaki"
$ws.Range("L12").Value = "Primera"
$ws.Range("M12").Value = 140
$ws.Range("N12").Value = 17000
$ws.Range("O12").Value = 18000
$ws.Range("P12").Value = 17429
$ws.Range("Q12").Value = '$/caja 18 kilos granel'
$ws.Range("R12").Value = "Provincia de Curicó"
$ws.Range("S12").Value = 968
$ws.Range("T12").Value = 18
